$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old row 3 (Nxph3/Nrxn1 -> MuSCs pair) entirely; row 2 becomes the
# only data row and the sheet dimension shrinks from A1:T3 to A1:T2.
$ws.Rows(3).Delete()

# Row 2 now describes the Nxph3/Nrxn1 pair signalling from ECs to MuSCs.
$ws.Range("A2").Value = "ECs"
$ws.Range("D2").Value = "MuSCs"

# Refresh row 2's numeric columns with the new TPM-derived values.
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.4198326666666667
$ws.Range("H2").Value = 1.259498
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.02430333333333333
$ws.Range("N2").Value = 0.07291
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 0.01020333324222222
$ws.Range("R2").Value = 0.09182999918000001
$ws.Range("S2").Value = 1
$ws.Range("T2").Value = 1
